# ISS update Dec 18
# Inserts a new most-recent row (computed on a later date) at the top of the
# data table, pushing the existing rows down by one, and re-stamps the
# Rt_reference_date (column C) time-of-day for every row to the new cutoff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (the first data row); this shifts
# rows 2:19 down to 3:20 and carries their formatting/styles with them.
$ws.Rows("2:2").Insert()

# The inserted row picks up formatting from the row above (the bold header),
# so re-apply the plain data-row formatting from the row below it instead.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New Rt_reference_date time-of-day offset applied to every data row.
$cOffset = 0.99930555555

# Newest data point (now row 2).
$ws.Cells.Item(2, 1).Value = 44160
$ws.Cells.Item(2, 2).Value = 44173
$ws.Cells.Item(2, 3).Value = 44166 + $cOffset
$ws.Cells.Item(2, 4).Value = 0.86
$ws.Cells.Item(2, 5).Value = 0.79
$ws.Cells.Item(2, 6).Value = 0.9399999999999999

# The rest of the series (rows 3:20), shifted down by one row from the
# previous layout; column C gets the refreshed time-of-day stamp.
$data = @(
    @(3,  44153, 44166, 44159, 0.82, 0.76, 0.91),
    @(4,  44146, 44159, 44152, 0.91, 0.79, 1.08),
    @(5,  44139, 44152, 44145, 1.08, 0.91, 1.25),
    @(6,  44132, 44145, 44138, 1.18, 0.9399999999999999, 1.49),
    @(7,  44126, 44139, 44132, 1.43, 1.08, 1.81),
    @(8,  44119, 44132, 44125, 1.72, 1.45, 1.83),
    @(9,  44112, 44125, 44118, 1.7, 1.49, 1.85),
    @(10, 44105, 44118, 44111, 1.5, 1.09, 1.75),
    @(11, 44098, 44111, 44104, 1.17, 1.03, 1.5),
    @(12, 44091, 44104, 44097, 1.06, 0.97, 1.16),
    @(13, 44084, 44097, 44090, 1.01, 0.88, 1.08),
    @(14, 44077, 44090, 44083, 0.95, 0.88, 1.05),
    @(15, 44070, 44083, 44076, 0.92, 0.79, 1.17),
    @(16, 44063, 44076, 44069, 1.14, 0.71, 1.53),
    @(17, 44056, 44069, 44062, 1.18, 0.86, 1.43),
    @(18, 44049, 44062, 44055, 0.75, 0.52, 1.24),
    @(19, 44042, 44055, 44048, 0.83, 0.67, 1.06),
    @(20, 44035, 44048, 44041, 0.96, 0.75, 1.2)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3] + $cOffset
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
